# Update countries & provincias Spain
#
# Refreshes the "paises.xlsx" COVID snapshot: the source numbers moved on
# (new totals for several countries), which shuffles a handful of
# countries' rank in the "Casos totales" (column B) descending sort, and
# bumps the "Datos actualizados" timestamp in the title row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column letter -> 1-based index, used to drive $ws.Cells.Item(row, col).
$colIndex = @{ A = 1; B = 2; C = 3; D = 4; E = 5; F = 6; G = 7; H = 8 }

# Rows whose metrics changed but keep their existing country name/rank.
$dataUpdates = @(
    @{ Row = 4;   D = 370977; E = 1127372 }
    @{ Row = 14;  B = 113321; C = 1293; D = 45900; E = 63965; G = 22;  H = 3456 }
    @{ Row = 28;  E = 1001;   G = 1;    H = 1893 }
    @{ Row = 55;  B = 8039;   C = 151;  D = 3715;  E = 4312 }
    @{ Row = 57;  D = 3777;   E = 3422 }
    @{ Row = 79;  B = 2812;   C = 98;   D = 1251;  E = 1531 }
    @{ Row = 80;  B = 2350;   C = 12;   D = 1596;  E = 614;  G = 4; H = 140 }
    @{ Row = 87;  B = 1898;   C = 40;   D = 1378;  E = 409;  G = 1; H = 111 }
    @{ Row = 118; B = 812;    C = 3;    D = 669;   E = 91 }
)

foreach ($u in $dataUpdates) {
    $row = $u.Row
    foreach ($col in 'B', 'C', 'D', 'E', 'F', 'G', 'H') {
        if ($u.ContainsKey($col)) {
            $ws.Cells.Item($row, $colIndex[$col]).Value = $u[$col]
        }
    }
}

# Rows where the refreshed totals also change the country sort order, so
# both the country name (col A) and its metrics are rewritten in place.
$rowRewrites = @(
    @{ Row = 125; A = 'Malta';                                B = 599; C = 15; D = 468; E = 125; H = 6 }
    @{ Row = 126; A = 'Haiti';                                B = 596;         D = 21;  E = 553; H = 22 }
    @{ Row = 136; A = 'Madagascar';                           B = 405; C = 34; D = 131; E = 272; H = 2 }
    @{ Row = 137; A = 'Etiopia';                                       C = 9;  D = 123; E = 270; H = 5 }
    @{ Row = 138; A = 'Estado de Palestina';                   B = 398;         D = 346; E = 50 }
    @{ Row = 196; A = 'Namibia';                                       C = 2;  D = 14;  E = 4 }
    @{ Row = 197; A = 'Fiyi';                                                  D = 15;  E = 3 }
    @{ Row = 198; A = 'Nueva Caledonia' }
    @{ Row = 199; A = 'Santa Lucia';                                           D = 18;           H = 0 }
    @{ Row = 200; A = 'Belice';                               B = 18;          D = 16;  E = 0;   H = 2 }
    @{ Row = 201; A = 'Islas Virgenes de los Estados Unidos'; B = 17;          D = 0;   E = 17 }
    @{ Row = 210; A = 'Groenlandia';                                           D = 11;           H = 0 }
    @{ Row = 211; A = 'Montserrat';                                            D = 10;           H = 1 }
    @{ Row = 214; A = 'Bonaire, San Eustaquio y Saba' }
    @{ Row = 215; A = 'Sahara Occidental' }
)

foreach ($u in $rowRewrites) {
    $row = $u.Row
    foreach ($col in 'A', 'B', 'C', 'D', 'E', 'F', 'G', 'H') {
        if ($u.ContainsKey($col)) {
            $ws.Cells.Item($row, $colIndex[$col]).Value = $u[$col]
        }
    }
}

# Bump the "last updated" timestamp in the title banner (row 1).
$ws.Range("A1").Value = "Datos actualizados a 21 de Mayo de 2020 a las 13:35"
